$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (losing exact formatting, e.g. "1.00" -> 1, "4.90" -> 4.9, "0.0000189" -> 1.89E-05).
# Pre-format them as Text so the literal string is preserved exactly, matching
# the inline-string cell content produced by the source report generator.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D12",
    "D16",
    "D20",
    "D21",
    "D22",
    "D23",
    "D25",
    "D27",
    "D29",
    "D31",
    "D38",
    "D39",
    "D41",
    "D47",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price (column D) and 1h-volume-change (column E) values
$ws.Range("D2").Value = '68.602.90'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.702.33'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '598.72'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '160.48'
$ws.Range("E6").Value = '  +2.91%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.545'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '2.701.26'
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '5.31'
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("D15").Value = '3.193.96'
$ws.Range("E15").Value = '  +2.26%  '
$ws.Range("D16").Value = '0.0000189'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '68.547.26'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '2.691.57'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").Value = '365.83'
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("D21").Value = '7.63'
$ws.Range("E21").Value = '  +4.03%  '
$ws.Range("D22").Value = '4.53'
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").Value = '  +2.56%  '
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").Value = '74.49'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '9.94'
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("D28").Value = '2.838.52'
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  +1.29%  '
$ws.Range("E30").Value = '  -6.86%  '
$ws.Range("D31").Value = '580.43'
$ws.Range("E31").Value = '  +4.58%  '
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("E34").Value = '  +5.61%  '
$ws.Range("E35").Value = '  +3.92%  '
$ws.Range("E36").Value = '  +6.38%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '161.77'
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").Value = '19.84'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  +2.16%  '
$ws.Range("D41").Value = '1.92'
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("E43").Value = '  +3.37%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  -4.79%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '157.93'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("E48").Value = '  +6.04%  '
$ws.Range("E49").Value = '  +5.20%  '
$ws.Range("E50").Value = '  +7.12%  '
$ws.Range("D51").Value = '22.08'
$ws.Range("E51").Value = '  +0.15%  '
